$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 776
$ws1.Range("F4").Value = 509
$ws1.Range("F5").Value = 257
$ws1.Range("F6").Value = 461
$ws1.Range("F7").Value = 1099
$ws1.Range("F11").Value = 99
$ws1.Range("F12").Value = 1088
$ws1.Range("F15").Value = 729
$ws1.Range("F16").Value = 785
$ws1.Range("F17").Value = 169
$ws1.Range("F19").Value = 50
$ws1.Range("F20").Value = 626
$ws1.Range("F21").Value = 119
$ws1.Range("F22").Value = 1689
$ws1.Range("F23").Value = 1930
$ws1.Range("F24").Value = 492
$ws1.Range("F26").Value = 1737
$ws1.Range("F28").Value = 2551
$ws1.Range("F29").Value = 458
$ws1.Range("F33").Value = 85
$ws1.Range("F34").Value = 89
$ws1.Range("F35").Value = 879
$ws1.Range("F36").Value = 1589
$ws1.Range("F37").Value = 267
$ws1.Range("F40").Value = 113
$ws1.Range("F41").Value = 98

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 9
$ws2.Range("F12").Value = 62

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 9
$ws4.Range("F5").Value = 776
$ws4.Range("F6").Value = 509
$ws4.Range("F7").Value = 257
$ws4.Range("F8").Value = 461
$ws4.Range("F9").Value = 1099
$ws4.Range("F13").Value = 99
$ws4.Range("F14").Value = 1088
$ws4.Range("F16").Value = 729
$ws4.Range("F17").Value = 785
$ws4.Range("F18").Value = 169
$ws4.Range("F24").Value = 50
$ws4.Range("F25").Value = 627
$ws4.Range("F26").Value = 119
$ws4.Range("F27").Value = 1689
$ws4.Range("F28").Value = 1930
$ws4.Range("F29").Value = 492
$ws4.Range("F32").Value = 2551
$ws4.Range("F33").Value = 458
$ws4.Range("F38").Value = 62
$ws4.Range("F41").Value = 85
$ws4.Range("F42").Value = 89
$ws4.Range("F43").Value = 879
$ws4.Range("F44").Value = 1589
$ws4.Range("F45").Value = 267
$ws4.Range("F47").Value = 113
$ws4.Range("F48").Value = 98
